$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 0.794582
$ws.Range("N2").Value = 2.383746
$ws.Range("O2").Value = 0.03449752952410986
$ws.Range("P2").Value = 0.03449752952410985
$ws.Range("Q2").Value = 0.1270539266606667
$ws.Range("R2").Value = 1.143485339946
$ws.Range("S2").Value = 0.0009134159079288777
$ws.Range("T2").Value = 0.0009134159079288777

$ws.Range("M3").Value = 20.604156
$ws.Range("N3").Value = 61.812468
$ws.Range("O3").Value = 0.8945489325574519
$ws.Range("P3").Value = 0.8945489325574517
$ws.Range("Q3").Value = 3.294611412452
$ws.Range("R3").Value = 29.651502712068
$ws.Range("S3").Value = 0.02368561565684628
$ws.Range("T3").Value = 0.02368561565684628

$ws.Range("M4").Value = 0.2871986666666667
$ws.Range("N4").Value = 0.8615959999999999
$ws.Range("O4").Value = 0.01246900191876775
$ws.Range("P4").Value = 0.01246900191876775
$ws.Range("Q4").Value = 0.04592316253288889
$ws.Range("R4").Value = 0.413308462796
$ws.Range("S4").Value = 0.0003301507344355856
$ws.Range("T4").Value = 0.0003301507344355855

$ws.Range("M5").Value = 1.149534666666667
$ws.Range("N5").Value = 3.448604
$ws.Range("O5").Value = 0.04990813547540859
$ws.Range("P5").Value = 0.04990813547540859
$ws.Range("Q5").Value = 0.1838109763782222
$ws.Range("R5").Value = 1.654298787404
$ws.Range("S5").Value = 0.001321453608625734
$ws.Range("T5").Value = 0.001321453608625734

$ws.Range("M6").Value = 0.1975403333333333
$ws.Range("N6").Value = 0.5926210000000001
$ws.Range("O6").Value = 0.008576400524262026
$ws.Range("P6").Value = 0.008576400524262026
$ws.Range("Q6").Value = 0.03158676514677778
$ws.Range("R6").Value = 0.284280886321
$ws.Range("S6").Value = 0.0002270835268408293
$ws.Range("T6").Value = 0.0002270835268408293

$ws.Range("M7").Value = 0.794582
$ws.Range("N7").Value = 2.383746
$ws.Range("O7").Value = 0.03449752952410986
$ws.Range("P7").Value = 0.03449752952410985
$ws.Range("Q7").Value = 4.671468354464
$ws.Range("R7").Value = 42.043215190176
$ws.Range("S7").Value = 0.03358411361618098
$ws.Range("T7").Value = 0.03358411361618097

$ws.Range("M8").Value = 20.604156
$ws.Range("N8").Value = 61.812468
$ws.Range("O8").Value = 0.8945489325574519
$ws.Range("P8").Value = 0.8945489325574517
$ws.Range("Q8").Value = 121.134964955712
$ws.Range("R8").Value = 1090.214684601408
$ws.Range("S8").Value = 0.8708633169006056
$ws.Range("T8").Value = 0.8708633169006055

$ws.Range("M9").Value = 0.2871986666666667
$ws.Range("N9").Value = 0.8615959999999999
$ws.Range("O9").Value = 0.01246900191876775
$ws.Range("P9").Value = 0.01246900191876775
$ws.Range("Q9").Value = 1.688484615530667
$ws.Range("R9").Value = 15.196361539776
$ws.Range("S9").Value = 0.01213885118433217
$ws.Range("T9").Value = 0.01213885118433216

$ws.Range("M10").Value = 1.149534666666667
$ws.Range("N10").Value = 3.448604
$ws.Range("O10").Value = 0.04990813547540859
$ws.Range("P10").Value = 0.04990813547540859
$ws.Range("Q10").Value = 6.758289034602667
$ws.Range("R10").Value = 60.824601311424
$ws.Range("S10").Value = 0.04858668186678286
$ws.Range("T10").Value = 0.04858668186678286

$ws.Range("M11").Value = 0.1975403333333333
$ws.Range("N11").Value = 0.5926210000000001
$ws.Range("O11").Value = 0.008576400524262026
$ws.Range("P11").Value = 0.008576400524262026
$ws.Range("Q11").Value = 1.161369645797333
$ws.Range("R11").Value = 10.452326812176
$ws.Range("S11").Value = 0.008349316997421197
$ws.Range("T11").Value = 0.008349316997421197
